$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new header columns (I, J) --------------------------------
# Copy the formatting of the last existing header cell (H1: bold, thin
# border, centered/top aligned) onto the two new header cells, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "MKL Time"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J1").Value = "CuBLAS Time"

# --- Existing trials (rows 2-30): stamp blank, typed cells in the two
#     new columns so the used range extends through column J ----------
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 9).Value = "'"
    $ws.Cells.Item($r, 9).Style = "Normal"
    $ws.Cells.Item($r, 10).Value = "'"
    $ws.Cells.Item($r, 10).Style = "Normal"
}

# --- New trials from the re-run with MKL / CuBLAS baselines -----------
$newTrials = @(
    @("Trial 30", 0, 0, 0, 0, 0, 0, 0.3082060813903809,  0.002991914749145508, 0.1581034660339355),
    @("Trial 31", 0, 0, 0, 0, 0, 0, 0.293776273727417,   0.002990961074829102, 0.1805171966552734),
    @("Trial 32", 0, 0, 0, 0, 0, 0, 0.3042478561401367,  0.003988981246948242, 0.1695470809936523),
    @("Trial 33", 0, 0, 0, 0, 0, 0, 0.3093512058258057,  0.004986286163330078, 0.1785213947296143),
    @("Trial 34", 0, 0, 0, 0, 0, 0, 0.3185737133026123,  0.0139617919921875,   0.183509349822998),
    @("Trial 35", 0, 0, 0, 0, 0, 0, 0.5135586261749268,  0.08676838874816895,  0.1924829483032227),
    @("Trial 36", 0, 0, 0, 0, 0, 0, 10.25245118141174,   4.482025146484375,    0.8058481216430664)
)

$startRow = 31
for ($i = 0; $i -lt $newTrials.Count; $i++) {
    $row = $startRow + $i
    $rowData = $newTrials[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
